$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("factory")
$ws.Rows.Item(13).Insert()
